$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.446.73"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "1.916.87"
$ws.Range("E3").Value = "  +2.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.97"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5111"
$ws.Range("E7").Value = "  +1.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3966"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09727"
$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.14"
$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.467"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.04"
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.920.07"
$ws.Range("E14").Value = "  +2.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.402"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.77"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06674"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("E20").Value = "  +4.44%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.259"
$ws.Range("E22").Value = "  +2.77%  "

$ws.Range("D23").Value = "28.488.73"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.322"
$ws.Range("E25").Value = "  +2.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.688"
$ws.Range("E26").Value = "  +6.81%  "

$ws.Range("D27").Value = "2.143.98"
$ws.Range("E27").Value = "  +2.75%  "

$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.40"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.11"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.107"
$ws.Range("E31").Value = "  +5.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1070"
$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.691"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.844"
$ws.Range("E35").Value = "  +5.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06705"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02440"
$ws.Range("E37").Value = "  +2.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.257"
$ws.Range("E38").Value = "  +4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2223"
$ws.Range("E39").Value = "  +2.04%  "

$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6427"
$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.222"
$ws.Range("E43").Value = "  +4.21%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.65"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6065"
$ws.Range("E46").Value = "  +1.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.783"
$ws.Range("E47").Value = "  +3.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.284"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.061"
$ws.Range("E49").Value = "  +3.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.70"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.197"
$ws.Range("E51").Value = "  +0.12%  "
